# The deck ships two theme parts:
#   theme1.xml  -> currently the stock "Office Theme" palette (only wired to
#                  the Notes Master, not visible anywhere in the main deck)
#   theme2.xml  -> currently the "Integral" palette (the one actually applied
#                  to the Slide Master / whole presentation)
#
# The authored edit swaps the two palettes: the Slide Master's theme becomes
# the plain "Office Theme" colors while the (unused) secondary theme becomes
# "Integral". We reproduce this by rewriting the 12 theme colors that PowerPoint
# exposes on the active design's ThemeColorScheme to the stock Office values -
# the same edit a user makes by picking the built-in "Office" color scheme
# from Design > Variants > Colors.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Index -> scheme slot (matches the Item() ordering PowerPoint exposes):
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
# Values are standard Office theme RGB colors, written as the 0xBBGGRR
# COM "RGB" integers PowerPoint's ColorFormat.RGB property expects
# (equivalent to calling the VBA RGB(r,g,b) helper).

$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
